$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.840.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.740.94"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5171"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2797"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.40"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06094"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.747.08"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07038"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.22"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6389"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.503"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.856.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006563"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.973.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.119"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.597"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.131"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.509"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.805"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.03"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08250"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.663"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.421"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04476"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.609"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9782"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6128"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.645"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01582"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.918"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.28"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3823"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7220"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.55%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.956"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05412"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.262"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1119"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.08"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.660"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("E51").Value = "  -0.80%  "
